# todolist.xlsx — add three new to-do rows and rewrite row 2.
#
# Target end state of the "Person" sheet (A1:B5):
#   Item | Priority
#   t    | low
#   ali  | low
#   5    | low
#   w    | low
#
# Row 2 (previously "ali" / "high") becomes "t" / "low", and three more
# rows are appended below it: "ali"/"low", "5"/"low", "w"/"low".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: overwrite the existing item/priority.
$ws.Range("A2").Value = "t"
$ws.Range("B2").Value = "low"

# Row 3: new item.
$ws.Range("A3").Value = "ali"
$ws.Range("B3").Value = "low"

# Row 4: new item "5" — force text storage (via a temporary Text number
# format) so it isn't auto-coerced into the number 5, then clear the
# formatting back off so the cell keeps the workbook's default style.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "5"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").Value = "low"

# Row 5: new item.
$ws.Range("A5").Value = "w"
$ws.Range("B5").Value = "low"
